# Fix property_category values that were incorrectly left as "land" for
# the building (建物) and car (汽車) sheets.

$wb = $excel.ActiveWorkbook

# 建物 (building) sheet: property_category column I, rows 2-3: land -> building
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2").Value = "building"
$wsBuilding.Range("I3").Value = "building"

# 汽車 (car) sheet: property_category column H, row 2: land -> car
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
